# Denormalize db schema: remove many-to-many relation table by adding
# "tag_ids" and "doc_ids" columns (comma-separated id lists) directly on
# the dataset table, replacing the old junction-table approach.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add the two new columns to the table (extends table ref + dimension,
# and keeps the tableColumns / worksheet columns in sync).
$colTagIds = $lo.ListColumns.Add()
$colDocIds = $lo.ListColumns.Add()

# Populate the new column data (order matters for shared-string layout:
# doc_ids data is written before tag_ids data).
$ws.Range("R1").Value = "doc_ids"
$ws.Range("R4").Value = "pdf_wiki, pdf_online"

$ws.Range("Q1").Value = "tag_ids"
$ws.Range("Q4").Value = "personal_data, sensible_data, sante, population"
$ws.Range("Q12").Value = "anonymous_data, population, culture"

# The long tag list on row 4 wraps within its cell.
$ws.Range("Q4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 16

# Column widths for the two new columns.
$ws.Columns.Item(17).ColumnWidth = 37
$ws.Columns.Item(18).ColumnWidth = 16
$ws.Columns.Item(19).ColumnWidth = 16

# Freeze header row + first column, with the view scrolled so column O is
# the left-most visible (unfrozen) column and E7 is the active cell.
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollColumn = 15
$ws.Range("E7").Select()
